$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("航天发展","航天发展","安妮股份"),
    @("龙洲股份","龙洲股份","航天发展"),
    @("航天动力","国脉科技","海王生物"),
    @("顺灏股份","航天科技","龙洲股份"),
    @("海王生物","海王生物","实达集团"),
    @("平潭发展","实达集团","平潭发展"),
    @("航天机电","东方财富","合富中国"),
    @("实达集团","顺灏股份","国脉科技"),
    @("国脉科技","赢时胜","顺灏股份"),
    @("合富中国","N摩尔-U","摩尔线程"),
    @("航天科技","航天机电","瑞康医药"),
    @("N摩尔-U","航天动力","国机重装"),
    @("东方财富","国机重装","航天动力"),
    @("乾照光电","西部材料","特发信息"),
    @("国机重装","中银证券","雷科防务"),
    @("西部材料","瑞康医药","航天机电"),
    @("中银证券","永鼎股份","海欣食品"),
    @("瑞康医药","特发信息","太阳电缆"),
    @("赢时胜","合富中国","安记食品"),
    @("特发信息","海欣食品","航天科技")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
